# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet named "2022-Q1" right before the "总计"
#    (totals) sheet, populated with the per-fund holdings for the new
#    quarter (same layout/formatting as the other quarterly sheets).
# 2. Prepend a new row to the "总计" sheet summarising the 2022-Q1
#    quarter (count of funds held + total market value), shifting the
#    previously existing rows down by one.

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# --- 1. Create the new "2022-Q1" sheet, positioned right after 2021-Q4
#        (i.e. right before 总计). ---------------------------------------
$newSheet = $wb.Worksheets.Add($null, $q4Sheet)
$newSheet.Name = "2022-Q1"

# Seed it from the 2021-Q4 layout (same headers + same visual style for
# the header row / index column), then overwrite every data cell below
# with the 2022-Q1 numbers. 2021-Q4 only has 5 fund rows, 2022-Q1 needs
# 6, so duplicate the last formatted row first to extend the range.
$q4Sheet.Range("A1:H6").Copy($newSheet.Range("A1"))
$newSheet.Range("A6:H6").Copy($newSheet.Range("A7"))

$rows = @(
    @(0, "090001", "大成价值增长混合", 18.65, 61.32, 4.74, 0.8840, 3),
    @(1, "506001", "万家科创板 2 年定期开放混合型证券投资基金", 12.84, 98.14, 4.31, 0.5534, 3),
    @(2, "160919", "大成产业升级股票(LOF)", 3.95, 87.76, 6.16, 0.2433, 1),
    @(3, "012210", "申万菱信智能汽车股票型证券投资基金A", 4.76, 82.52, 3.88, 0.1847, 9),
    @(4, "012051", "申万菱信乐道三年持有期混合型证券投资基金", 3.38, 81.64, 4.12, 0.1393, 7),
    @(5, "012211", "申万菱信智能汽车股票型证券投资基金C", 1.40, 82.52, 3.88, 0.0543, 9)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $newSheet.Cells.Item($r, 1).Value = $data[0]

    # Fund code: force text so leading zeros ("090001", "012210", ...)
    # survive instead of being parsed as a number.
    $codeCell = $newSheet.Cells.Item($r, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $data[1]

    $newSheet.Cells.Item($r, 3).Value = $data[2]
    $newSheet.Cells.Item($r, 4).Value = $data[3]
    $newSheet.Cells.Item($r, 5).Value = $data[4]
    $newSheet.Cells.Item($r, 6).Value = $data[5]
    $newSheet.Cells.Item($r, 7).Value = $data[6]
    $newSheet.Cells.Item($r, 8).Value = $data[7]
}

# --- 2. Insert the 2022-Q1 summary row at the top of "总计" (below the
#        header), pushing the existing rows down. ------------------------
# Re-resolve the totals sheet AFTER the insert above: sheet references
# here track by position, so grabbing it earlier would now (post-insert)
# point at the freshly added "2022-Q1" sheet instead.
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 2.06

# Renumber the leading index column (0,1,2,...) for every data row, since
# Insert() shifts cell content down but leaves the old index values as-is.
for ($r = 2; $r -le 6; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

Write-Output "done"
